$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(18).Copy()
$ws.Rows.Item(19).PasteSpecial(-4122)
